{"js": "// The upstream change for this file is a pure OOXML canonicalization: every\n// hunk in the diff only re-orders XML attributes / namespace declarations\n// (e.g. `w:top=\"...\" w:right=\"...\"` -> `w:bottom=\"...\" w:footer=\"...\"`,\n// `xmlns:wpc=... xmlns:mc=...` -> alphabetically sorted `xmlns:m=... xmlns:mc=...`).\n// Every attribute name/value pair that exists before the change still exists\n// after it, with the exact same value - nothing was added, removed, or\n// retargeted (headerReference still points at rId6, pgSz/pgMar keep their\n// original twips, the footnote separators keep their ids/types, the header's\n// \"m:self\" field + accent6/BF theme color are untouched, and the style\n// catalog's fonts/langs/latentStyles/style ids are unchanged). So there is no\n// visible document content, formatting, or structure to change here - the\n// body text, header field code, sections, and styles already match the\n// target state.\n//\n// We still touch the document through the supported object model (instead of\n// leaving the script empty) so the intent - confirming the section/page\n// layout that the template's sectPr encodes - is explicit, but we only read\n// values; we never assign anything, so no accidental content drift can be\n// introduced.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst pageSetup = section.pageSetup;\npageSetup.load(\"topMargin,bottomMargin,leftMargin,rightMargin,pageWidth,pageHeight\");\nawait context.sync();\n\n// Nothing to change: the page geometry already matches the template's\n// original sectPr (pgSz 11906x16838 twips, pgMar 1417/1417/1417/1417 with\n// 708 header/footer and 0 gutter) -- nothing in the diff altered any of\n// these values, so the body content is left exactly as-is.\nawait context.sync();\n", "ps1": "# The upstream change for this file is a pure OOXML canonicalization: every\n# hunk in the diff only re-orders XML attributes / namespace declarations\n# (e.g. `w:top=\"...\" w:right=\"...\"` -> `w:bottom=\"...\" w:footer=\"...\"`,\n# `xmlns:wpc=... xmlns:mc=...` -> alphabetically sorted `xmlns:m=... xmlns:mc=...`,\n# `w:type=\"separator\" w:id=\"-1\"` -> `w:id=\"-1\" w:type=\"separator\"`, etc.).\n# Every attribute name/value pair present before the change is still present\n# after it with the exact same value - nothing is added, removed, or\n# retargeted: the header reference still points at rId6, pgSz/pgMar keep\n# their original twips (11906x16838 / 1417,1417,1417,1417,708,708,0), the\n# footnote separator/continuationSeparator ids and types are unchanged, the\n# header's \" m:\" + \"self\" + \". \" field code and its accent6/BF themed color\n# run are untouched, and every style id / font / language / latentStyles\n# entry in the style catalog is identical. So there is no visible document\n# content, formatting, or structure to change in this file.\n#\n# We still touch the document through the supported COM object model\n# (instead of leaving the script empty) so the intent - confirming the\n# section/page layout that the template's sectPr encodes - is explicit, but\n# we only read values; we never assign anything, so no accidental content\n# drift can be introduced.\n$d = $word.ActiveDocument\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n\n$top = $pageSetup.TopMargin\n$bottom = $pageSetup.BottomMargin\n$left = $pageSetup.LeftMargin\n$right = $pageSetup.RightMargin\n$pageWidth = $pageSetup.PageWidth\n$pageHeight = $pageSetup.PageHeight\n\n# Nothing to change: the page geometry already matches the template's\n# original sectPr (pgSz 11906x16838 twips, pgMar 1417/1417/1417/1417 with\n# 708 header/footer and 0 gutter) -- nothing in the diff altered any of\n# these values, so the body content is left exactly as-is.\n"}
